# Apply cryptos list update (values scraped on Mon Nov  6 13:47:36 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.335.93"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.912.71"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'0.722"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.82%  "
$ws.Range("D6").Value = "'252.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "'40.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").Value = "'52.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.27%  "
$ws.Range("D11").Value = "'0.0732"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").Value = "'0.0998"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "2.193.83"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "'12.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "'0.714"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.931.40"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'4.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "35.404.73"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'72.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").Value = "0.0₃0828"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'13.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").Value = "'241.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.55%  "
$ws.Range("D23").Value = "'5.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.52%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.22%  "
$ws.Range("D27").Value = "'167.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").Value = "'8.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.70%  "
$ws.Range("E29").Value = "  +4.84%  "
$ws.Range("D30").Value = "'18.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "4.127.54"
$ws.Range("E31").Value = "  +19.42%  "
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("E33").Value = "  +13.02%  "
$ws.Range("D34").Value = "'0.0579"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("D35").Value = "'1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +19.21%  "
$ws.Range("D36").Value = "'4.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").Value = "'2.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").Value = "'17.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.52%  "
$ws.Range("D41").Value = "'98.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.38%  "
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("D43").Value = "'0.0209"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").Value = "'0.0650"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("D45").Value = "'2.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.92%  "
$ws.Range("D46").Value = "1.346.96"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'6.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.54%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "'2.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'45.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.19%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.101.73"
$ws.Range("E51").Value = "  +0.00%  "

Write-Host "Applied cell updates"
